$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet from "1" to "ონი"
$ws.Name = "ონი"

# Remove the census-note row (old row 2: "(მოსახლეობის აღწერის შედეგებით)").
# This shifts the old empty row 3 up to become the new (empty) row 2.
$ws.Rows("2").Delete()

# Remove the 1989 and 2002 columns, keeping only the 2014 column (old column D,
# now becomes column B).
$ws.Columns("B:C").Delete()

# Put the selection on A2 to match the saved view state.
$ws.Range("A2").Select() | Out-Null
